# Quick agenda update - 18:13:34,71
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - camera issue resolved: clear the "Teste" follow-up column, status -> Concluido
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "Concluido"

# Row 3 - Roberto / Localiza ticket description + observation updated, status -> Pendente
$ws.Range("D3").Value = "Arrumar sensores externos, estão sem funcionar."
$ws.Range("E3").Value = "Roberto esteve no local mas aparentemente só deu inicio."
$ws.Range("G3").Value = "Pendente"

# Row 4 - new entry: Pedro (itauna) / Casa da Kenia
$ws.Range("A4").Value = "Pedro (itaúna)"
$ws.Range("B4").Value = "'2693"
$ws.Range("C4").Value = "Casa da Kênia"
$ws.Range("D4").Value = "Colocar central via internet."
$ws.Range("E4").Value = "Foi instalado o módulo, agora funciona via internet."
$ws.Range("F4").Value = "Foi colocado o módulo no local."
$ws.Range("G4").Value = "Concluido"

# Row 5 - new entry: Pedro (itauna) / Guia CWK
$ws.Range("A5").Value = "Pedro (itaúna)"
$ws.Range("B5").Value = "2138"
$ws.Range("C5").Value = "Guia CWK"
$ws.Range("D5").Value = "Local sem comunicação de alarmes."
$ws.Range("G5").Value = "Pendente"

# Row 6 - new entry: Pedro (itauna) / Rodonaves
$ws.Range("A6").Value = "Pedro (itaúna)"
$ws.Range("B6").Value = "2565"
$ws.Range("C6").Value = "Rodonaves"
$ws.Range("D6").Value = "Reparo e acesso em câmeras não monitoradas."
$ws.Range("G6").Value = "Pendente"

# Move selection to G6 (scrolled back to A1/no fixed topLeftCell)
$ws.Range("G6").Select()
